# Fruta / hortaliza, semanal
# Swap the weekly data between row 2 <-> row 5 and row 3 <-> row 6
# (columns D, M, N, O, P, R, S carry the differing values; the rest of
# each row's columns are identical between the swapped pairs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($row1, $row2, $col) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($col in $cols) {
    Swap-Cell 2 5 $col
}

foreach ($col in $cols) {
    Swap-Cell 3 6 $col
}
